# The edit reorders the 4 data rows (rows 2-5) of the active sheet: each
# row's full contents move up by one position (row3->row2, row4->row3,
# row5->row4), and the original row 2 wraps around to become row 5. In
# other words, a cyclic rotation of rows [2,3,4,5] -> [3,4,5,2].
#
# Row.Copy(destination) in this host only overwrites cells that are
# non-blank in the source range, leaving stale values behind in any
# destination cell whose source counterpart is blank. So every
# destination range is explicitly cleared immediately before the copy
# to guarantee a faithful move (including cells that should end up
# blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol = "AY"
$tempRow = 100

function RowRange($row) {
    return $ws.Range($firstCol + $row + ":" + $lastCol + $row)
}

# 1) Stash original row 2 in a scratch row well below the used range.
(RowRange $tempRow).ClearContents()
(RowRange 2).Copy((RowRange $tempRow))

# 2) Shift rows 3,4,5 up into 2,3,4.
(RowRange 2).ClearContents()
(RowRange 3).Copy((RowRange 2))

(RowRange 3).ClearContents()
(RowRange 4).Copy((RowRange 3))

(RowRange 4).ClearContents()
(RowRange 5).Copy((RowRange 4))

# 3) Move the stashed original row 2 into row 5.
(RowRange 5).ClearContents()
(RowRange $tempRow).Copy((RowRange 5))

# 4) Clean up the scratch row.
(RowRange $tempRow).ClearContents()
